$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# The six "optimal hyperparameter / metrics" sheets (arbolcc, bosquecc,
# knncc, arbolpp, bosquepp, knnpp) each had a one-row results table in
# A1:F2:
#   A: hyperparameter label (suffixed with the sheet name)
#   B: "Arreglo aleatorio óptimo <sheet>" (random-state search result)
#   C: MAE, D: MSE, E: RMSE, F: "R-cuadrado <sheet>"
#
# The new model run drops column B entirely (random-state column no longer
# reported), renames the remaining headers (strip the sheet-name suffix
# from column A, rename "R-cuadrado" -> "R2"), and reports fresh metric
# values - leaving a tight A1:E2 table (A: hyperparam, B: MAE, C: MSE,
# D: RMSE, E: R2).
#
# Other sheets (e.g. ResltNumericas) hold formulas like "=arbolcc!B2" /
# "=arbolcc!F2" pointing at specific columns of these tables, and those
# formulas are untouched by this edit - so we overwrite cell values in
# place (A-E) and wipe column F with a full Clear() (contents + format)
# rather than doing a structural column delete/insert, which would shift
# and break those external references.
# ---------------------------------------------------------------------------

function Update-ResultSheet {
    param(
        [string]$SheetName,
        [string]$ColATitle,
        $ColAValue,
        $MAE,
        $MSE,
        $RMSE,
        $R2
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Row 1 headers.
    $ws.Range("A1").Value = $ColATitle
    $ws.Range("B1").Value = "MAE " + $SheetName
    $ws.Range("C1").Value = "MSE " + $SheetName
    $ws.Range("D1").Value = "RMSE " + $SheetName
    $ws.Range("E1").Value = "R2 " + $SheetName

    # Row 2 values.
    $ws.Range("A2").Value = $ColAValue
    $ws.Range("B2").Value = $MAE
    $ws.Range("C2").Value = $MSE
    $ws.Range("D2").Value = $RMSE
    $ws.Range("E2").Value = $R2

    # Drop the now-unused column F (old "R-cuadrado ..." column) entirely -
    # content + formatting - so the used range/dimension shrinks to A1:E2.
    $ws.Range("F1:F2").Clear()
}

Update-ResultSheet "arbolcc" "Profundidad óptima" `
    3 0.5774729620742622 0.6466338346268923 `
    0.8041354578843618 0.3648838465338654

Update-ResultSheet "bosquecc" "Estimador óptimo" `
    115 0.4293526781499747 0.3504596228933595 `
    0.5919963031078483 0.6557826768132963

Update-ResultSheet "knncc" "K óptimo" `
    3 0.5326407001766574 0.4800962164278529 `
    0.6928897577738128 0.5284551380655178

Update-ResultSheet "arbolpp" "Profundidad óptima" `
    6 0.775569221353989 1.827375618630338 `
    1.351804578565385 -1.115941009930692

Update-ResultSheet "bosquepp" "Estimador óptimo" `
    110 0.6276056324687056 0.5735812803590546 `
    0.7573514906297172 0.3358430848771266

Update-ResultSheet "knnpp" "K óptimo" `
    5 0.7557783359119404 0.8315282410723833 `
    0.9118817034420547 0.03716308335154872
